$d = $word.ActiveDocument

# --- Step 1: remove the old "Meta description" paragraph (the 2nd paragraph,
#     right under the title heading) ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- Step 2: insert a new bold "Play Eye of the Storm..." paragraph right
#     before the final paragraph ("Create a feature image..."). First add a
#     plain paragraph break, then fill that new, empty paragraph with the
#     exact run XML (a leading empty <w:r/> plus the bold text run) via
#     Range.InsertXML - this also overwrites any inherited paragraph
#     formatting, so the new paragraph ends up with no <w:pPr> at all, just
#     like the rest of the document. ---
$cnt = $d.Paragraphs.Count
$secondLastPara = $d.Paragraphs.Item($cnt - 1)
[void]$secondLastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($cnt)
$newRange = $newPara.Range

$runXmlPkg = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Eye of the Storm Free Online Slot - Review 2021</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$newRange.InsertXML($runXmlPkg)

# --- Step 3: replace the text of the final paragraph (keeps its original
#     italic run formatting, since Find/Replace only swaps the text) ---
$oldText = "Create a feature image for " + [char]34 + "Eye of the Storm" + [char]34 + ": In your design, incorporate a cartoon-style image of a happy Maya warrior. The warrior should be wearing glasses and surrounded by Egyptian-themed symbols such as the Eye of Ra, pyramids, and hieroglyphs. Place the warrior at the center of the image, with the Eye of Ra symbol behind him. Ensure that the color scheme is vibrant and eye-catching, with a focus on gold, black, and beige tones. Add some lightning bolts in the background and depict the warrior as if he's ready to take on the Eye of the Storm and grab all the riches waiting to be discovered!"
$newText = "Read our review of Eye of the Storm online slot game. Play for free with its stunning graphics, Egyptian theme and Wild and Scatter for big wins."

[void]$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                               $true, 1, $false, $newText, 2)
